$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.698.92"
$ws.Range("E2").Value = "  -1.45%  "

$ws.Range("D3").Value = "2.458.79"
$ws.Range("E3").Value = "  -2.05%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.75"
$ws.Range("E5").Value = "  -1.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.35"
$ws.Range("E6").Value = "  -2.10%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -1.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.152"
$ws.Range("E9").Value = "  -6.27%  "

$ws.Range("E10").Value = "  -1.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.343"
$ws.Range("E11").Value = "  -3.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.82"
$ws.Range("E12").Value = "  -2.51%  "

$ws.Range("D13").Value = "2.911.39"
$ws.Range("E13").Value = "  -2.02%  "

$ws.Range("D14").Value = "68.549.12"
$ws.Range("E14").Value = "  -1.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  -3.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.73"
$ws.Range("E16").Value = "  -4.61%  "

$ws.Range("D17").Value = "2.479.38"
$ws.Range("E17").Value = "  -1.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.04"
$ws.Range("E18").Value = "  -1.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "345.19"
$ws.Range("E19").Value = "  -1.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.22"
$ws.Range("E20").Value = "  -3.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.83"
$ws.Range("E21").Value = "  -2.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.90"
$ws.Range("E22").Value = "  -3.46%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.22"
$ws.Range("E24").Value = "  -2.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.77"
$ws.Range("E25").Value = "  -4.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.04"
$ws.Range("E26").Value = "  +4.41%  "

$ws.Range("D27").Value = "2.582.88"
$ws.Range("E27").Value = "  -2.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.27"
$ws.Range("E28").Value = "  -6.51%  "

$ws.Range("D29").Value = "0.0₃0845"
$ws.Range("E29").Value = "  -5.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.34"
$ws.Range("E30").Value = "  -6.71%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  -2.85%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "437.24"
$ws.Range("E32").Value = "  -4.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("E34").Value = "  -2.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.07"
$ws.Range("E35").Value = "  +104.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.97"
$ws.Range("E36").Value = "  -1.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.02"
$ws.Range("E37").Value = "  -0.25%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  -5.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.94"
$ws.Range("E40").Value = "  -3.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.308"
$ws.Range("E41").Value = "  -3.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.53"
$ws.Range("E42").Value = "  -3.66%  "

$ws.Range("E43").Value = "  -3.67%  "

$ws.Range("E44").Value = "  +1.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.11"
$ws.Range("E45").Value = "  -4.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.88"
$ws.Range("E46").Value = "  -4.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("E47").Value = "  -2.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.491"
$ws.Range("E48").Value = "  -5.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0721"
$ws.Range("E49").Value = "  -1.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.566"
$ws.Range("E50").Value = "  -2.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0916"
$ws.Range("E51").Value = "  -1.49%  "

